$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "report_comment"
$ws.Range("E7").Value = "Recent points demonstrate special-cause improvement.  Congratulations and carry on!"
$ws.Range("E2").Value = "This is a comment about the attendances metric.  This text is added via 'report_config.xlsx'"

$ws.Columns.Item(5).ColumnWidth = 72

$ws.Range("E3").Select()
